# Korrigiere Projekt- und WBS-ID-Referenzen in der Main.controller.js
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAPUI5 Export")

# Row 2: projectId (col A) must hold the bare project id; wbsId (col B)
# keeps the full "project.wbs" reference it already had.
$ws.Range("A2").Value = "10.30.00002"

# Row 3: same correction for the second project.
$ws.Range("A3").Value = "10.20.00019"

# New row 4: additional WBS element for project 10.20.00019, reusing the
# same planned/baseline dates as row 3 and flagged as "Long Project".
$ws.Range("C3:F3").Copy()
$ws.Range("C4:F4").PasteSpecial(-4122)   # xlPasteFormats - keep the date style

$ws.Range("A4").Value = "10.20.00019"
$ws.Range("B4").Value = "10.20.00019.1010103"
$ws.Range("C4").Value = 45719
$ws.Range("D4").Value = 45780
$ws.Range("E4").Value = 45719
$ws.Range("F4").Value = 45780
$ws.Range("H4").Value = "Long Project"

$ws.Range("D12").Select()
